$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "96-48=48"
$t.Cell(1, 2).Range.Text = "86-83=3"
$t.Cell(1, 3).Range.Text = "99-65=34"
$t.Cell(1, 4).Range.Text = "64-11=53"
$t.Cell(1, 5).Range.Text = "56+10=66"

$t.Cell(2, 1).Range.Text = "84-41=43"
$t.Cell(2, 2).Range.Text = "93-40=53"
$t.Cell(2, 3).Range.Text = "96-2=94"
$t.Cell(2, 4).Range.Text = "99-66=33"
$t.Cell(2, 5).Range.Text = "23+55=78"

$t.Cell(3, 1).Range.Text = "38+38=76"
$t.Cell(3, 2).Range.Text = "49-12=37"
$t.Cell(3, 3).Range.Text = "93-60=33"
$t.Cell(3, 4).Range.Text = "60-32=28"
$t.Cell(3, 5).Range.Text = "57-33=24"

$t.Cell(4, 1).Range.Text = "79+7=86"
$t.Cell(4, 2).Range.Text = "67-39=28"
$t.Cell(4, 3).Range.Text = "91-42=49"
$t.Cell(4, 4).Range.Text = "89-29=60"
$t.Cell(4, 5).Range.Text = "35-23=12"

$t.Cell(5, 1).Range.Text = "7+67=74"
$t.Cell(5, 2).Range.Text = "28-9=19"
$t.Cell(5, 3).Range.Text = "96-2=94"
$t.Cell(5, 4).Range.Text = "74-23=51"
$t.Cell(5, 5).Range.Text = "67-40=27"

$t.Cell(6, 1).Range.Text = "27-8=19"
$t.Cell(6, 2).Range.Text = "38-5=33"
$t.Cell(6, 3).Range.Text = "93-60=33"
$t.Cell(6, 4).Range.Text = "65+9=74"
$t.Cell(6, 5).Range.Text = "30-11=19"

$t.Cell(7, 1).Range.Text = "30+47=77"
$t.Cell(7, 2).Range.Text = "57+19=76"
$t.Cell(7, 3).Range.Text = "81-50=31"
$t.Cell(7, 4).Range.Text = "86-49=37"
$t.Cell(7, 5).Range.Text = "61-43=18"

$t.Cell(8, 1).Range.Text = "67+14=81"
$t.Cell(8, 2).Range.Text = "44+22=66"
$t.Cell(8, 3).Range.Text = "50+8=58"
$t.Cell(8, 4).Range.Text = "27+32=59"
$t.Cell(8, 5).Range.Text = "10+2=12"

$t.Cell(9, 1).Range.Text = "56+21=77"
$t.Cell(9, 2).Range.Text = "75-72=3"
$t.Cell(9, 3).Range.Text = "81-47=34"
$t.Cell(9, 4).Range.Text = "96-59=37"
$t.Cell(9, 5).Range.Text = "64-57=7"

$t.Cell(10, 1).Range.Text = "93-45=48"
$t.Cell(10, 2).Range.Text = "64-3=61"
$t.Cell(10, 3).Range.Text = "17+52=69"
$t.Cell(10, 4).Range.Text = "79-20=59"
$t.Cell(10, 5).Range.Text = "14+43=57"

$t.Cell(11, 1).Range.Text = "0+21=21"
$t.Cell(11, 2).Range.Text = "61-10=51"
$t.Cell(11, 3).Range.Text = "53-6=47"
$t.Cell(11, 4).Range.Text = "61+20=81"
$t.Cell(11, 5).Range.Text = "45-13=32"

$t.Cell(12, 1).Range.Text = "89-47=42"
$t.Cell(12, 2).Range.Text = "1+92=93"
$t.Cell(12, 3).Range.Text = "2+12=14"
$t.Cell(12, 4).Range.Text = "63+30=93"
$t.Cell(12, 5).Range.Text = "21+76=97"

$t.Cell(13, 1).Range.Text = "39+9=48"
$t.Cell(13, 2).Range.Text = "6+25=31"
$t.Cell(13, 3).Range.Text = "6+27=33"
$t.Cell(13, 4).Range.Text = "44+35=79"
$t.Cell(13, 5).Range.Text = "79-56=23"

$t.Cell(14, 1).Range.Text = "98-60=38"
$t.Cell(14, 2).Range.Text = "86-61=25"
$t.Cell(14, 3).Range.Text = "4+38=42"
$t.Cell(14, 4).Range.Text = "94-18=76"
$t.Cell(14, 5).Range.Text = "37-33=4"

$t.Cell(15, 1).Range.Text = "66-21=45"
$t.Cell(15, 2).Range.Text = "5+45=50"
$t.Cell(15, 3).Range.Text = "60-6=54"
$t.Cell(15, 4).Range.Text = "36+29=65"
$t.Cell(15, 5).Range.Text = "1+10=11"

$t.Cell(16, 1).Range.Text = "38+61=99"
$t.Cell(16, 2).Range.Text = "48+7=55"
$t.Cell(16, 3).Range.Text = "39-4=35"
$t.Cell(16, 4).Range.Text = "27+28=55"
$t.Cell(16, 5).Range.Text = "92-4=88"

$t.Cell(17, 1).Range.Text = "68-61=7"
$t.Cell(17, 2).Range.Text = "47-45=2"
$t.Cell(17, 3).Range.Text = "90-34=56"
$t.Cell(17, 4).Range.Text = "36+9=45"
$t.Cell(17, 5).Range.Text = "22+26=48"

$t.Cell(18, 1).Range.Text = "68-7=61"
$t.Cell(18, 2).Range.Text = "0+36=36"
$t.Cell(18, 3).Range.Text = "29+2=31"
$t.Cell(18, 4).Range.Text = "69-66=3"
$t.Cell(18, 5).Range.Text = "2+9=11"

$t.Cell(19, 1).Range.Text = "72-33=39"
$t.Cell(19, 2).Range.Text = "91-19=72"
$t.Cell(19, 3).Range.Text = "64-52=12"
$t.Cell(19, 4).Range.Text = "28-26=2"
$t.Cell(19, 5).Range.Text = "90-27=63"

$t.Cell(20, 1).Range.Text = "82+14=96"
$t.Cell(20, 2).Range.Text = "77+11=88"
$t.Cell(20, 3).Range.Text = "38+21=59"
$t.Cell(20, 4).Range.Text = "99-43=56"
$t.Cell(20, 5).Range.Text = "13+75=88"

